# Updates cryptocurrency price/volume figures (and two row swaps) to refresh
# the latest "cryptos" snapshot, per the GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.137.71"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.649.57"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D5").Value = "'218.54"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.2625"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "'0.06310"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'20.39"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "'0.07652"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "'4.582"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "1.637.38"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "1.876.10"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "'0.5582"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "0.0₅8129"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "'65.17"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "26.089.40"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'4.597"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").Value = "'194.54"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'10.48"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'145.31"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "'7.202"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'15.85"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").Value = "'0.05478"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "'1.271"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "'3.452"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "'3.330"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "'1.561"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "'2.782"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.412"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "'0.9429"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'0.5620"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "'0.01574"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.747"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "1.028.82"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'0.8202"
$ws.Range("E43").Value = "  -3.17%  "
$ws.Range("D44").Value = "'100.62"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").Value = "1.787.46"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +6.91%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'57.31"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "'0.4323"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "'7.932"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "'0.05125"
$ws.Range("E51").Value = "  -3.56%  "
